$d = $word.ActiveDocument

# Locate the empty "List Paragraph" that immediately follows
# "Created a Repository folder and created two repositories for the two
# models created in the previous step." -- this is the bullet that the
# commit turns into the "Created Dummy Mocks (...)" bullet, right before
# the new "Added service configurations..." bullet is inserted after it.
$targetPara = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -eq "`r" -and $p.Range.ParagraphFormat.Style.NameLocal -eq "List Paragraph") {
        $targetPara = $p
        break
    }
}

if ($null -eq $targetPara) {
    throw "Could not locate the target empty list paragraph"
}

# Build the exact OOXML for the two paragraphs:
#  1) the previously-empty paragraph (keeps its original pPr / paraId)
#     now filled with the "Created Dummy Mocks (...)" runs + proofErr
#     spell-check markers around the two mock file names, and
#  2) a brand-new "List Paragraph" bullet right after it describing the
#     service configuration / statup.cs change.
$xml = @'
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:body><w:p w14:paraId="14808345" w14:textId="77777777" w:rsidR="00391EF0" w:rsidRPr="002C44BF" w:rsidRDefault="00391EF0" w:rsidP="002C44BF"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>Created Dummy Mocks</w:t></w:r><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>(</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>MockCategoryRepository.cs</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> and </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>MockItemRepository.cs</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>)</w:t></w:r><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> in a newly created Mocks folder to provide data to </w:t></w:r><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>the</w:t></w:r><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> application using the created Interfaces</w:t></w:r><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>/Repositories.</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">Added service configurations and implementations for the above created interfaces and dummy mocks in the Configure Services section of the </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>statup.cs</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> file.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@

$targetPara.Range.InsertXML($xml)

Write-Host "Inserted 'Created Dummy Mocks...' and 'Added service configurations...' bullets."
